$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47-72 down to 48-73.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new weekly record.
$ws.Range("A47").Value = 5
$ws.Range("B47").Value = "Macroferia Regional de Talca"
$ws.Range("C47").Value = "Maule"
$ws.Range("D47").Value = (Get-Date -Year 2022 -Month 11 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E47").Value = 7
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100101
$ws.Range("H47").Value = "Berries"
$ws.Range("I47").Value = 100101001
$ws.Range("J47").Value = "Arándano (blue)"
$ws.Range("K47").Value = "Sin especificar"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 35
$ws.Range("N47").Value = 3500
$ws.Range("O47").Value = 3500
$ws.Range("P47").Value = 3500
$ws.Range("Q47").Value = "$/bandeja 2 kilos"
$ws.Range("R47").Value = "Provincia de Curicó"
$ws.Range("S47").Value = 1750
$ws.Range("T47").Value = 2
